$d = $word.ActiveDocument

# Update the date line (unique text, safe to use Find/Replace)
$d.Content.Find.Execute("2025-11-01 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-11-02 Sunday", 2)

# Update the division problems table. The table has 20 rows x 5 columns,
# with only every 4th row (1, 5, 9, 13, 17) containing data.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "810÷4=202, 2"
$t.Cell(1, 2).Range.Text  = "869÷9=96, 5"
$t.Cell(1, 3).Range.Text  = "152÷4=38, 0"
$t.Cell(1, 4).Range.Text  = "584÷7=83, 3"
$t.Cell(1, 5).Range.Text  = "480÷9=53, 3"

$t.Cell(5, 1).Range.Text  = "710÷5=142, 0"
$t.Cell(5, 2).Range.Text  = "168÷6=28, 0"
$t.Cell(5, 3).Range.Text  = "329÷3=109, 2"
$t.Cell(5, 4).Range.Text  = "304÷4=76, 0"
$t.Cell(5, 5).Range.Text  = "280÷2=140, 0"

$t.Cell(9, 1).Range.Text  = "205÷4=51, 1"
$t.Cell(9, 2).Range.Text  = "705÷7=100, 5"
$t.Cell(9, 3).Range.Text  = "223÷9=24, 7"
$t.Cell(9, 4).Range.Text  = "937÷3=312, 1"
$t.Cell(9, 5).Range.Text  = "188÷4=47, 0"

$t.Cell(13, 1).Range.Text = "103÷4=25, 3"
$t.Cell(13, 2).Range.Text = "315÷3=105, 0"
$t.Cell(13, 3).Range.Text = "227÷8=28, 3"
$t.Cell(13, 4).Range.Text = "451÷3=150, 1"
$t.Cell(13, 5).Range.Text = "318÷4=79, 2"

$t.Cell(17, 1).Range.Text = "678÷2=339, 0"
$t.Cell(17, 2).Range.Text = "733÷8=91, 5"
$t.Cell(17, 3).Range.Text = "842÷9=93, 5"
$t.Cell(17, 4).Range.Text = "268÷6=44, 4"
$t.Cell(17, 5).Range.Text = "330÷5=66, 0"
